$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new value looks like a plain decimal number
# need an explicit Text number format first, otherwise Excel COM auto-converts
# the assigned string into a floating point number (losing the original
# "37.820.89"-style / trailing-zero text formatting used throughout this sheet).
$ws.Cells.Item(2, 4).Value = "37.802.28"
$ws.Cells.Item(2, 5).Value = "  +1.17%  "
$ws.Cells.Item(3, 4).Value = "2.084.25"
$ws.Cells.Item(3, 5).Value = "  +0.88%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "232.43"
$ws.Cells.Item(5, 5).Value = "  -0.73%  "
$ws.Cells.Item(6, 5).Value = "  -0.39%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "57.24"
$ws.Cells.Item(8, 5).Value = "  +0.48%  "
$ws.Cells.Item(9, 5).Value = "  +1.44%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0778"
$ws.Cells.Item(10, 5).Value = "  +2.14%  "
$ws.Cells.Item(11, 5).Value = "  +2.76%  "
$ws.Cells.Item(12, 4).Value = "2.383.06"
$ws.Cells.Item(12, 5).Value = "  +0.58%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "14.36"
$ws.Cells.Item(13, 5).Value = "  -1.68%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "21.02"
$ws.Cells.Item(14, 5).Value = "  +1.80%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.758"
$ws.Cells.Item(15, 5).Value = "  -2.44%  "
$ws.Cells.Item(16, 5).Value = "  +1.95%  "
$ws.Cells.Item(17, 4).Value = "2.083.59"
$ws.Cells.Item(17, 5).Value = "  +0.89%  "
$ws.Cells.Item(18, 4).Value = "37.732.16"
$ws.Cells.Item(18, 5).Value = "  +1.14%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.14"
$ws.Cells.Item(19, 5).Value = "  -2.30%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "70.83"
$ws.Cells.Item(20, 5).Value = "  +1.98%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0819"
$ws.Cells.Item(21, 5).Value = "  +1.16%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "227.87"
$ws.Cells.Item(22, 5).Value = "  +0.64%  "
$ws.Cells.Item(23, 5).Value = "  -0.05%  "
$ws.Cells.Item(24, 5).Value = "  -2.03%  "
$ws.Cells.Item(25, 5).Value = "  -1.19%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "169.81"
$ws.Cells.Item(26, 5).Value = "  +1.97%  "
$ws.Cells.Item(27, 5).Value = "  +9.78%  "
$ws.Cells.Item(28, 5).Value = "  +1.37%  "
$ws.Cells.Item(29, 5).Value = "  -0.19%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "19.43"
$ws.Cells.Item(30, 5).Value = "  +1.87%  "
$ws.Cells.Item(31, 5).Value = "  +0.74%  "
$ws.Cells.Item(32, 5).Value = "  +2.68%  "
$ws.Cells.Item(34, 5).Value = "  +0.36%  "
$ws.Cells.Item(35, 5).Value = "  +0.64%  "
$ws.Cells.Item(36, 5).Value = "  +3.71%  "
$ws.Cells.Item(37, 5).Value = "  +4.07%  "
$ws.Cells.Item(38, 5).Value = "  +0.07%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.40"
$ws.Cells.Item(39, 5).Value = "  -4.76%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0994"
$ws.Cells.Item(40, 5).Value = "  +5.89%  "
$ws.Cells.Item(41, 5).Value = "  -0.81%  "
$ws.Cells.Item(42, 5).Value = "  +0.84%  "
$ws.Cells.Item(43, 5).Value = "  +0.28%  "
$ws.Cells.Item(44, 4).Value = "1.452.47"
$ws.Cells.Item(44, 5).Value = "  -0.92%  "
$ws.Cells.Item(45, 5).Value = "  -1.05%  "
$ws.Cells.Item(46, 5).Value = "  +3.01%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.05"
$ws.Cells.Item(47, 5).Value = "  -7.81%  "
$ws.Cells.Item(48, 5).Value = "  +3.73%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.36"
$ws.Cells.Item(49, 5).Value = "  +3.03%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.99"
$ws.Cells.Item(50, 5).Value = "  +1.28%  "
$ws.Cells.Item(51, 4).Value = "2.278.02"
$ws.Cells.Item(51, 5).Value = "  +0.88%  "
